# Auto-generated script applying scheduled market-data refresh to Tonberry_Profits sheets.
# For each affected Leve row, currentAveragePrice/-NQ/-HQ, LevePriceNQ/HQ and
# LeveProfitNQ/HQ columns (H:N) are updated to the latest market-board snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 12076.186
$ws.Range("I18").Value = 7219.7856
$ws.Range("J18").Value = 17306.154
$ws.Range("K18").Value = 7219.7856
$ws.Range("L18").Value = 17306.154
$ws.Range("M18").Value = -6935.7856
$ws.Range("N18").Value = -17874.154

$ws.Range("H64").Value = 2853.3333
$ws.Range("I64").Value = 2750
$ws.Range("J64").Value = 2982.5
$ws.Range("K64").Value = 2750
$ws.Range("L64").Value = 2982.5
$ws.Range("M64").Value = -2502
$ws.Range("N64").Value = -3478.5

$ws.Range("H67").Value = 2853.3333
$ws.Range("I67").Value = 2750
$ws.Range("J67").Value = 2982.5
$ws.Range("K67").Value = 2750
$ws.Range("L67").Value = 2982.5
$ws.Range("M67").Value = -1892
$ws.Range("N67").Value = -4698.5

$ws.Range("H100").Value = 2181
$ws.Range("I100").Value = 1066.8572
$ws.Range("K100").Value = 1066.8572
$ws.Range("M100").Value = -525.8571999999999

$ws.Range("H125").Value = 2478.6667
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 2478.6667
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 22308.0003
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -27228.0003

$ws.Range("H127").Value = 919.1429000000001
$ws.Range("I127").Value = 839.1667
$ws.Range("J127").Value = 1399
$ws.Range("K127").Value = 2517.5001
$ws.Range("L127").Value = 4197
$ws.Range("M127").Value = 2442.4999
$ws.Range("N127").Value = -14117

$ws.Range("H135").Value = 517.5625
$ws.Range("I135").Value = 465.07693
$ws.Range("K135").Value = 4185.69237
$ws.Range("M135").Value = -1650.69237

$ws.Range("H137").Value = 1905.0625
$ws.Range("I137").Value = 1125.6666
$ws.Range("K137").Value = 3376.9998
$ws.Range("M137").Value = -826.9998000000001

$ws.Range("H138").Value = 2204.4211
$ws.Range("I138").Value = 1798.125
$ws.Range("J138").Value = 2499.9092
$ws.Range("K138").Value = 5394.375
$ws.Range("L138").Value = 7499.7276
$ws.Range("M138").Value = -254.375
$ws.Range("N138").Value = -17779.7276


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1787.1837
$ws.Range("I32").Value = 1310.5796
$ws.Range("K32").Value = 1310.5796
$ws.Range("M32").Value = -1023.5796

$ws.Range("H74").Value = 1339.0952
$ws.Range("I74").Value = 600.5714
$ws.Range("J74").Value = 2816.1428
$ws.Range("K74").Value = 600.5714
$ws.Range("L74").Value = 2816.1428
$ws.Range("M74").Value = 273.4286
$ws.Range("N74").Value = -4564.1428

$ws.Range("H77").Value = 1339.0952
$ws.Range("I77").Value = 600.5714
$ws.Range("J77").Value = 2816.1428
$ws.Range("K77").Value = 3002.857
$ws.Range("L77").Value = 14080.714
$ws.Range("M77").Value = 1365.143
$ws.Range("N77").Value = -22816.714

$ws.Range("H101").Value = 71000
$ws.Range("J101").Value = 71000
$ws.Range("L101").Value = 71000
$ws.Range("N101").Value = -77490

$ws.Range("H122").Value = 49080.125
$ws.Range("I122").Value = 64773.5
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 194320.5
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -191870.5
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 2096.875
$ws.Range("I132").Value = 1970.0333
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 5910.0999
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -3380.0999
$ws.Range("N132").Value = -17058.5


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2611
$ws.Range("I105").Value = 2588.1177
$ws.Range("K105").Value = 2588.1177
$ws.Range("M105").Value = -841.1176999999998

$ws.Range("H134").Value = 5600.3076
$ws.Range("I134").Value = 5733.6665
$ws.Range("K134").Value = 17200.9995
$ws.Range("M134").Value = -14665.9995


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 650
$ws.Range("I105").Value = 500
$ws.Range("J105").Value = 800
$ws.Range("K105").Value = 500
$ws.Range("L105").Value = 800
$ws.Range("M105").Value = 1247
$ws.Range("N105").Value = -4294

$ws.Range("H132").Value = 2523.8823
$ws.Range("I132").Value = 1810.9286
$ws.Range("K132").Value = 5432.7858
$ws.Range("M132").Value = -2902.7858

$ws.Range("H134").Value = 1392.6552
$ws.Range("I134").Value = 1116.8462
$ws.Range("K134").Value = 3350.5386
$ws.Range("M134").Value = -815.5385999999999


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 1997.25
$ws.Range("I70").Value = 999.5
$ws.Range("J70").Value = 2995
$ws.Range("K70").Value = 2998.5
$ws.Range("L70").Value = 8985
$ws.Range("M70").Value = -2683.5
$ws.Range("N70").Value = -9615

$ws.Range("H73").Value = 1997.25
$ws.Range("I73").Value = 999.5
$ws.Range("J73").Value = 2995
$ws.Range("K73").Value = 2998.5
$ws.Range("L73").Value = 8985
$ws.Range("M73").Value = -1906.5
$ws.Range("N73").Value = -11169

$ws.Range("H75").Value = 32591.834
$ws.Range("I75").Value = 599
$ws.Range("J75").Value = 38990.4
$ws.Range("K75").Value = 1797
$ws.Range("L75").Value = 116971.2
$ws.Range("M75").Value = -799
$ws.Range("N75").Value = -118967.2

$ws.Range("H78").Value = 32591.834
$ws.Range("I78").Value = 599
$ws.Range("J78").Value = 38990.4
$ws.Range("K78").Value = 5391
$ws.Range("L78").Value = 350913.6
$ws.Range("M78").Value = -399
$ws.Range("N78").Value = -360897.6

$ws.Range("H81").Value = 19876962
$ws.Range("I81").Value = 2280
$ws.Range("J81").Value = 28910908
$ws.Range("K81").Value = 6840
$ws.Range("L81").Value = 86732724
$ws.Range("M81").Value = -5717
$ws.Range("N81").Value = -86734970

$ws.Range("H84").Value = 19876962
$ws.Range("I84").Value = 2280
$ws.Range("J84").Value = 28910908
$ws.Range("K84").Value = 20520
$ws.Range("L84").Value = 260198172
$ws.Range("M84").Value = -14904
$ws.Range("N84").Value = -260209404

$ws.Range("H87").Value = 12176
$ws.Range("I87").Value = 802.6667
$ws.Range("J87").Value = 19000
$ws.Range("K87").Value = 2408.0001
$ws.Range("L87").Value = 57000
$ws.Range("M87").Value = -1160.0001
$ws.Range("N87").Value = -59496

$ws.Range("H90").Value = 12176
$ws.Range("I90").Value = 802.6667
$ws.Range("J90").Value = 19000
$ws.Range("K90").Value = 7224.0003
$ws.Range("L90").Value = 171000
$ws.Range("M90").Value = -984.0002999999997
$ws.Range("N90").Value = -183480

$ws.Range("H108").Value = 2162.8333
$ws.Range("I108").Value = 825.6667
$ws.Range("J108").Value = 3500
$ws.Range("K108").Value = 2477.0001
$ws.Range("L108").Value = 10500
$ws.Range("M108").Value = 402.9998999999998
$ws.Range("N108").Value = -16260

$ws.Range("H116").Value = 2442.6667
$ws.Range("J116").Value = 3000
$ws.Range("L116").Value = 9000
$ws.Range("N116").Value = -15884

$ws.Range("H131").Value = 10888200
$ws.Range("J131").Value = 20865.098
$ws.Range("L131").Value = 62595.29400000001
$ws.Range("N131").Value = -72675.29400000001

$ws.Range("H132").Value = 984.1905
$ws.Range("I132").Value = 800
$ws.Range("J132").Value = 1014.8889
$ws.Range("K132").Value = 7200
$ws.Range("L132").Value = 9134.000100000001
$ws.Range("M132").Value = -4670
$ws.Range("N132").Value = -14194.0001


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 25000
$ws.Range("J86").Value = 25000
$ws.Range("L86").Value = 25000
$ws.Range("N86").Value = -27372

$ws.Range("H89").Value = 25000
$ws.Range("J89").Value = 25000
$ws.Range("L89").Value = 75000
$ws.Range("N89").Value = -86856

$ws.Range("H102").Value = 4178.0586
$ws.Range("I102").Value = 4201.8667
$ws.Range("J102").Value = 3999.5
$ws.Range("K102").Value = 4201.8667
$ws.Range("L102").Value = 3999.5
$ws.Range("M102").Value = -2579.8667
$ws.Range("N102").Value = -7243.5

$ws.Range("H113").Value = 1875
$ws.Range("I113").Value = 1733.3334
$ws.Range("J113").Value = 1960
$ws.Range("K113").Value = 1733.3334
$ws.Range("L113").Value = 1960
$ws.Range("M113").Value = 436.6666
$ws.Range("N113").Value = -6300

$ws.Range("H122").Value = 2345.389
$ws.Range("J122").Value = 2517.9
$ws.Range("L122").Value = 7553.700000000001
$ws.Range("N122").Value = -12453.7

$ws.Range("H132").Value = 1203845.9
$ws.Range("I132").Value = 1604156.6
$ws.Range("J132").Value = 2913.375
$ws.Range("K132").Value = 4812469.800000001
$ws.Range("L132").Value = 8740.125
$ws.Range("M132").Value = -4809939.800000001
$ws.Range("N132").Value = -13800.125


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 29999
$ws.Range("J47").Value = 29999
$ws.Range("L47").Value = 29999
$ws.Range("N47").Value = -30979

$ws.Range("H52").Value = 29999
$ws.Range("J52").Value = 29999
$ws.Range("L52").Value = 29999
$ws.Range("N52").Value = -30465

$ws.Range("H100").Value = 2484.9167
$ws.Range("I100").Value = 2258.0908
$ws.Range("K100").Value = 2258.0908
$ws.Range("M100").Value = -1717.0908

$ws.Range("H132").Value = 3795.0833
$ws.Range("I132").Value = 1299.5
$ws.Range("J132").Value = 4021.9546
$ws.Range("K132").Value = 3898.5
$ws.Range("L132").Value = 12065.8638
$ws.Range("M132").Value = -1368.5
$ws.Range("N132").Value = -17125.8638

$ws.Range("H136").Value = 3347.6562
$ws.Range("I136").Value = 2185.625
$ws.Range("K136").Value = 6556.875
$ws.Range("M136").Value = -4006.875


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 45630.71
$ws.Range("I122").Value = 56252.12
$ws.Range("K122").Value = 168756.36
$ws.Range("M122").Value = -166306.36

$ws.Range("H132").Value = 1411.0233
$ws.Range("I132").Value = 1252.8667
$ws.Range("J132").Value = 1776
$ws.Range("K132").Value = 3758.6001
$ws.Range("L132").Value = 5328
$ws.Range("M132").Value = -1228.6001
$ws.Range("N132").Value = -10388

$ws.Range("H136").Value = 11823536
$ws.Range("I136").Value = 19160048
$ws.Range("J136").Value = 3599.1667
$ws.Range("K136").Value = 57480144
$ws.Range("L136").Value = 10797.5001
$ws.Range("M136").Value = -57477594
$ws.Range("N136").Value = -15897.5001

